# Generate Report for Archive
# - Update "Ready for handoff" status to "In Translation" everywhere it appears
#   (Overview sheet's zh-cn/de-de status columns, and each language sheet's
#   "Status" column).
# - Narrow the corresponding "Status" columns to their new (smaller) width.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet (zh-cn / de-de status columns E & F, rows 2-3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$usedOverview = $wsOverview.UsedRange
for ($r = 1; $r -le $usedOverview.Rows.Count; $r++) {
    foreach ($colLetter in @("E", "F")) {
        $cell = $wsOverview.Range("$colLetter$r")
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# Narrow columns E and F on the Overview sheet.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-language sheets (zh-cn, de-de): "Status" column ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange

    # Find the "Status" column by reading the header row.
    $statusCol = 0
    for ($c = 1; $c -le $used.Columns.Count; $c++) {
        if ($ws.Cells.Item(1, $c).Value2 -eq "Status") {
            $statusCol = $c
        }
    }

    if ($statusCol -gt 0) {
        for ($r = 1; $r -le $used.Rows.Count; $r++) {
            $cell = $ws.Cells.Item($r, $statusCol)
            if ($cell.Value2 -eq $oldStatus) {
                $cell.Value = $newStatus
            }
        }

        # Narrow the Status column to match the report's new layout.
        $ws.Columns.Item($statusCol).ColumnWidth = 12.5
    }
}
